$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename
$ws.Range("O1").Value = "F1 train"

# Column O value updates for rows 2-10
$ws.Range("O2").Value = 0.9620253164556962
$ws.Range("O3").Value = 1
$ws.Range("O4").Value = 0.987012987012987
$ws.Range("O5").Value = 0.810126582278481
$ws.Range("O6").Value = 0.6823529411764706
$ws.Range("O7").Value = 1
$ws.Range("O8").Value = 0.987012987012987
$ws.Range("O9").Value = 1
$ws.Range("O10").Value = 0.821917808219178

# Row 11 updates
$ws.Range("C11").Value = "{'activation': 'tanh', 'alpha': 0.0001, 'hidden_layer_sizes': (64, 32), 'learning_rate': 'constant'}"
$ws.Range("E11").Value = 7
$ws.Range("F11").Value = 5
$ws.Range("G11").Value = 5
$ws.Range("H11").Value = 3
$ws.Range("J11").Value = 0.6363636363636364
$ws.Range("K11").Value = 0.7
$ws.Range("L11").Value = 0.5833333333333334
$ws.Range("M11").Value = 0.5
$ws.Range("N11").Value = 0.7
$ws.Range("O11").Value = 0.6666666666666666

# Column O value updates for rows 12-15
$ws.Range("O12").Value = 1
$ws.Range("O13").Value = 1
$ws.Range("O14").Value = 1
$ws.Range("O15").Value = 0.7619047619047619

# Row 16 updates
$ws.Range("C16").Value = "{'activation': 'tanh', 'alpha': 0.0001, 'hidden_layer_sizes': (64, 32), 'learning_rate': 'constant'}"
$ws.Range("E16").Value = 8
$ws.Range("H16").Value = 2
$ws.Range("I16").Value = 0.75
$ws.Range("J16").Value = 0.7619047619047619
$ws.Range("K16").Value = 0.8
$ws.Range("L16").Value = 0.7272727272727273
$ws.Range("N16").Value = 0.8
$ws.Range("O16").Value = 0.7848101265822784
